$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the Dictionary worksheet entirely - it is no longer needed
$wb.Worksheets("Dictionary").Delete() | Out-Null

# Add two more classified statements (sports related) to Sentiments_Arabic
$ws1 = $wb.Worksheets("Sentiments_Arabic")
$ws1.Range("A11").Value = "برشلونة تخسر مباراتها أمام البايرن"
$ws1.Range("B11").Value = "رياضة"
$ws1.Range("A12").Value = "خسارة فريق السودان أمام السنغال"
$ws1.Range("B12").Value = "رياضة"

# Add seven new rows of sentiment samples to Sentiments_Analysis,
# copying the formatting of the last existing data row
$ws2 = $wb.Worksheets("Sentiments_Analysis")
$ws2.Range("A8:B8").Copy() | Out-Null
$ws2.Range("A9:B15").PasteSpecial(-4122) | Out-Null

$ws2.Range("A9").Value = "You are stupid"
$ws2.Range("B9").Value = "Bad"

$ws2.Range("A10").Value = "Never, I will not"
$ws2.Range("B10").Value = "Bad"

$ws2.Range("A11").Value = "Excellent dear"
$ws2.Range("B11").Value = "Good"

$ws2.Range("A12").Value = "Good "
$ws2.Range("B12").Value = "Good"

$ws2.Range("A13").Value = "Bad"
$ws2.Range("B13").Value = "Bad"

$ws2.Range("A14").Value = "Nice"
$ws2.Range("B14").Value = "Good"

$ws2.Range("A15").Value = "Never"
$ws2.Range("B15").Value = "Bad"

# Update selections left behind by editing, and make Sentiments_Analysis
# the active tab/sheet
$ws1.Select() | Out-Null
$ws1.Range("B13").Select() | Out-Null

$ws2.Select() | Out-Null
$ws2.Range("A16").Select() | Out-Null
